# NewLeadCreation changes: add an "Execute" status column, a "Rating"/"Source"
# pair of columns, and three new lead rows to Sheet1; tweak Sheet2's view.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 structural changes -------------------------------------------

# The existing mailto hyperlink on the Email column doesn't get its cell
# reference updated by a later column insert, so drop it now and re-add it
# (and the three new ones) once the grid is in its final shape.
$ws1.Hyperlinks.Delete()

# Insert a new first column -> "Execute" (shifts FirstName..Country right).
$ws1.Columns("A:A").Insert()

# Insert two more columns before the (now) Country column -> "Rating" and
# "Source" (shifts Country from F to I).
$ws1.Columns("G:H").Insert()

$ws1.Columns("H:H").ColumnWidth = 14.95

# --- Header row -------------------------------------------------------
$ws1.Range("A1").Value = "Execute"
$ws1.Range("G1").Value = "Rating"
$ws1.Range("H1").Value = "Source"

# --- Row 2 (existing lead, now with Execute/Rating/Source filled in) ---
$ws1.Range("A2").Value = "No"
$ws1.Range("G2").Value = "Warm"
$ws1.Range("H2").Value = "Advertisement"

# --- Row 3 (new lead: John Doe) ----------------------------------------
$ws1.Range("A3").Value = "Yes"
$ws1.Range("B3").Value = "John"
$ws1.Range("C3").Value = "Doe"
$ws1.Range("D3").Value = "Manager"
$ws1.Range("E3").Value = "Accenture"
$ws1.Range("F3").Value = "johndoe@accenture.com"
$ws1.Range("G3").Value = "Hot"
$ws1.Range("H3").Value = "Customer Event"
$ws1.Range("I3").Value = "England"
$ws1.Range("A3:I3").Borders.LineStyle = 1

# --- Row 4 (new lead: Hillary Doe) --------------------------------------
$ws1.Range("A4").Value = "Yes"
$ws1.Range("B4").Value = "Hillary"
$ws1.Range("C4").Value = "Doe"
$ws1.Range("D4").Value = "VP"
$ws1.Range("E4").Value = "Wipro"
$ws1.Range("F4").Value = "hdoe@wipro.com"
$ws1.Range("G4").Value = "Cold"
$ws1.Range("H4").Value = "Partner"
$ws1.Range("I4").Value = "France"
$ws1.Range("A4:I4").Borders.LineStyle = 1

# --- Row 5 (new lead: Caleb Roy) ----------------------------------------
$ws1.Range("A5").Value = "No"
$ws1.Range("B5").Value = "Caleb"
$ws1.Range("C5").Value = "Roy"
$ws1.Range("D5").Value = "Director"
$ws1.Range("E5").Value = "Apple"
$ws1.Range("F5").Value = "caleb@apple.ca"
$ws1.Range("G5").Value = "Warm"
$ws1.Range("H5").Value = "Employee Referral"
$ws1.Range("I5").Value = "Australia"
$ws1.Range("A5:E5").Borders(7).LineStyle = 1
$ws1.Range("A5:E5").Borders(10).LineStyle = 1
$ws1.Range("G5:I5").Borders(7).LineStyle = 1
$ws1.Range("G5:I5").Borders(10).LineStyle = 1

# --- Header "Execute" cell gets a distinct fill -------------------------
$ws1.Range("A1").Interior.ThemeColor = 6
$ws1.Range("A1").Borders.LineStyle = 1

# --- Hyperlinks (re-created after the structural edits above) ----------
$ws1.Hyperlinks.Add($ws1.Range("F2"), "mailto:josephjefries@gmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F3"), "mailto:johndoe@accenture.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F4"), "mailto:hdoe@wipro.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F5"), "mailto:caleb@apple.ca") | Out-Null
$ws1.Range("F3").Borders.LineStyle = 1
$ws1.Range("F4").Borders.LineStyle = 1

# --- Selection on Sheet1 -------------------------------------------------
$ws1.Range("A5").Select()

# --- Sheet2 view tweak (selection moves; keep Sheet1 the active tab) ---
$ws2.Activate()
$ws2.Range("M8").Select()
$ws1.Activate()
